$wb = $excel.ActiveWorkbook

$names = @(
    "summ05871647",
    "summ06088680",
    "summ06318854",
    "summ06548205",
    "summ06767754",
    "summ07019143",
    "summ07251678",
    "summ07468429",
    "summ07698778"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $names[$i - 1]
}
